$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "CubeA"

# Fix tiny floating point update on L15 (rounding difference from re-export)
$ws.Range("L15").Value = 0.8954914636872461

# Add new row 16 of averaged-intensity data (Gaussian Quadrature scheme output)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.7618749660278448
$ws.Range("D16").Value = 1.342963186989
$ws.Range("E16").Value = 0.9519747158388456
$ws.Range("F16").Value = 1.054547556704468
$ws.Range("G16").Value = 0.7618749660278448
$ws.Range("H16").Value = 1.342963186989
$ws.Range("I16").Value = 0.8963520568295001
$ws.Range("J16").Value = 1.059311935314379
$ws.Range("K16").Value = 0.8977606274731419
$ws.Range("L16").Value = 1.201449708837466
$ws.Range("M16").Value = 0.7618749660278448
$ws.Range("N16").Value = 1.147468951413923
$ws.Range("O16").Value = 1.02784010639004
$ws.Range("P16").Value = 1.02077934425183

# Match the style used by A2:A15 (bold, centered, top-aligned, bordered)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
